$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.860.62"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.041.38"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").Value = "2.343.63"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.780"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("E16").Value = "  -2.67%  "

$ws.Range("D17").Value = "2.042.69"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "37.836.12"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.59"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.12%  "

$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("E24").Value = "  -1.81%  "

$ws.Range("E25").Value = "  +2.65%  "

$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("E28").Value = "  -3.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("E37").Value = "  +2.33%  "

$ws.Range("E38").Value = "  +4.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.32%  "

$ws.Range("D41").Value = "1.534.77"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("E44").Value = "  -1.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "

$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("E48").Value = "  -0.48%  "

$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "2.231.82"
$ws.Range("E51").Value = "  -0.71%  "
